$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns whose values are swapped between row 3 and row 4
$cols = @("A", "B", "D", "E", "F", "G", "H", "S")

foreach ($col in $cols) {
    $cell3 = $ws.Range($col + "3")
    $cell4 = $ws.Range($col + "4")
    $v3 = $cell3.Value2
    $v4 = $cell4.Value2
    $cell3.Value2 = $v4
    $cell4.Value2 = $v3
}
